$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text looks numeric need NumberFormat forced to Text
# temporarily so Excel does not coerce/trim the string (e.g. "13.70" -> 13.7),
# then the style is reset to Normal so no residual formatting is left behind.
$numericPriceRows = @(5,6,20,21,22,23,25,27,28,30,33,34,38,39,40,41,42,46,48)
foreach ($r in $numericPriceRows) {
    $ws.Range("D" + $r).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.036.42"
$ws.Range("D3").Value = "2.998.11"
$ws.Range("D5").Value = "560.33"
$ws.Range("D6").Value = "136.73"
$ws.Range("D9").Value = "2.992.41"
$ws.Range("D16").Value = "3.487.65"
$ws.Range("D18").Value = "2.996.03"
$ws.Range("D19").Value = "59.050.41"
$ws.Range("D20").Value = "428.92"
$ws.Range("D21").Value = "13.70"
$ws.Range("D22").Value = "0.721"
$ws.Range("D23").Value = "7.11"
$ws.Range("D25").Value = "80.65"
$ws.Range("D27").Value = "1.00"
$ws.Range("D28").Value = "2.18"
$ws.Range("D30").Value = "7.81"
$ws.Range("D33").Value = "0.0994"
$ws.Range("D34").Value = "0.996"
$ws.Range("D36").Value = "0.0₃0761"
$ws.Range("D38").Value = "48.75"
$ws.Range("D39").Value = "8.66"
$ws.Range("D40").Value = "2.74"
$ws.Range("D41").Value = "399.20"
$ws.Range("D42").Value = "0.0352"
$ws.Range("D43").Value = "2.756.49"
$ws.Range("D46").Value = "35.59"
$ws.Range("D48").Value = "123.39"

foreach ($r in $numericPriceRows) {
    $ws.Range("D" + $r).Style = "Normal"
}

$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("E11").Value = "  +7.41%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +6.89%  "
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("E28").Value = "  +9.18%  "
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E34").Value = "  +5.80%  "
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("E36").Value = "  +8.76%  "
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  +5.66%  "
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("E46").Value = "  +23.77%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("E51").Value = "  -1.00%  "
